{"js": "// The worksheet date heading moves forward one day (Tue 2025-04-22 -> Wed 2025-04-23),\n// and every arithmetic problem in the 20x5 practice table is replaced with a new\n// expression. Only the w:t text content changes; run/paragraph formatting is untouched.\n\nconst NEW_DATE = \"2025-04-23 Wednesday\";\n\n// New text for every cell, in row-major order (20 rows x 5 columns), matching the\n// table's existing layout exactly.\nconst NEW_CELL_VALUES = [\n  [\"66+25=\", \"25+27=\", \"85-63=\", \"28+22=\", \"84-70=\"],\n  [\"13+78=\", \"48+3=\", \"98-57=\", \"41+48=\", \"11+37=\"],\n  [\"52-1=\", \"71-46=\", \"50+1=\", \"15+15=\", \"45+54=\"],\n  [\"17-10=\", \"96-28=\", \"8+39=\", \"77-65=\", \"83+2=\"],\n  [\"26-17=\", \"62+0=\", \"12+41=\", \"20+37=\", \"37+33=\"],\n  [\"54+0=\", \"69-51=\", \"4+73=\", \"5-5=\", \"35+8=\"],\n  [\"46-38=\", \"27+64=\", \"39+35=\", \"3+25=\", \"6+90=\"],\n  [\"32+0=\", \"10+37=\", \"38+33=\", \"24-7=\", \"78-35=\"],\n  [\"30+56=\", \"76+8=\", \"13-5=\", \"70+14=\", \"51-47=\"],\n  [\"40-17=\", \"65-7=\", \"97-53=\", \"73-17=\", \"94-76=\"],\n  [\"10-2=\", \"6+11=\", \"66-22=\", \"33+12=\", \"94-21=\"],\n  [\"65-42=\", \"73+6=\", \"91-1=\", \"42+53=\", \"52-13=\"],\n  [\"62-27=\", \"69-31=\", \"14+84=\", \"71-43=\", \"66-23=\"],\n  [\"0+92=\", \"22+34=\", \"52-21=\", \"6+11=\", \"17+11=\"],\n  [\"48+21=\", \"95-77=\", \"42-25=\", \"98-33=\", \"93-33=\"],\n  [\"2+69=\", \"37+27=\", \"98-63=\", \"57+29=\", \"96-44=\"],\n  [\"93-19=\", \"38-33=\", \"28-26=\", \"61-47=\", \"41-18=\"],\n  [\"78-44=\", \"84-81=\", \"70-48=\", \"81-9=\", \"65-63=\"],\n  [\"85-63=\", \"89-37=\", \"86-67=\", \"63+22=\", \"70-11=\"],\n  [\"23+1=\", \"71-19=\", \"58-29=\", \"97-1=\", \"35-19=\"],\n];\n\nconst body = context.document.body;\n\n// 1. Update the date heading paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(NEW_DATE, Word.InsertLocation.replace);\n\n// 2. Update every cell of the practice table in a single shot via the `values`\n// property, which rewrites each cell's text while preserving its existing\n// paragraph/run formatting (font, size, alignment, etc.).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nif (table.rowCount !== NEW_CELL_VALUES.length) {\n  throw new Error(\n    `Unexpected row count: table has ${table.rowCount}, expected ${NEW_CELL_VALUES.length}`\n  );\n}\n\ntable.values = NEW_CELL_VALUES;\n\nawait context.sync();\n", "ps1": "# Worksheet date moves forward one day; every arithmetic problem in the\n# 20x5 practice table is replaced with a new expression. Only the cell text\n# changes -- existing run/paragraph formatting is left untouched.\n\n$newDate = '2025-04-23 Wednesday'\n\n# New text for every cell, in row-major order (20 rows x 5 columns).\n$newCellValues = @(\n    @('66+25=', '25+27=', '85-63=', '28+22=', '84-70='),\n    @('13+78=', '48+3=', '98-57=', '41+48=', '11+37='),\n    @('52-1=', '71-46=', '50+1=', '15+15=', '45+54='),\n    @('17-10=', '96-28=', '8+39=', '77-65=', '83+2='),\n    @('26-17=', '62+0=', '12+41=', '20+37=', '37+33='),\n    @('54+0=', '69-51=', '4+73=', '5-5=', '35+8='),\n    @('46-38=', '27+64=', '39+35=', '3+25=', '6+90='),\n    @('32+0=', '10+37=', '38+33=', '24-7=', '78-35='),\n    @('30+56=', '76+8=', '13-5=', '70+14=', '51-47='),\n    @('40-17=', '65-7=', '97-53=', '73-17=', '94-76='),\n    @('10-2=', '6+11=', '66-22=', '33+12=', '94-21='),\n    @('65-42=', '73+6=', '91-1=', '42+53=', '52-13='),\n    @('62-27=', '69-31=', '14+84=', '71-43=', '66-23='),\n    @('0+92=', '22+34=', '52-21=', '6+11=', '17+11='),\n    @('48+21=', '95-77=', '42-25=', '98-33=', '93-33='),\n    @('2+69=', '37+27=', '98-63=', '57+29=', '96-44='),\n    @('93-19=', '38-33=', '28-26=', '61-47=', '41-18='),\n    @('78-44=', '84-81=', '70-48=', '81-9=', '65-63='),\n    @('85-63=', '89-37=', '86-67=', '63+22=', '70-11='),\n    @('23+1=', '71-19=', '58-29=', '97-1=', '35-19='),\n)\n\n$d = $word.ActiveDocument\n\n# 1. Update the date heading (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = $newDate\n\n# 2. Update every cell of the practice table in place.\n$t = $d.Tables.Item(1)\nfor ($r = 0; $r -lt $newCellValues.Count; $r++) {\n    $row = $newCellValues[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n\nWrite-Output \"Updated date + $($newCellValues.Count * $newCellValues[0].Count) cells.\"\n"}
